# "Removendo e reorganizando Lista de Riscos"
#
# The risk row that used to be row 6 (Adriano/Waltson - "Feedback" risk,
# with the shared H6:H8 magnitude formula and its own one-off formatting)
# is removed from the "Riscos" sheet. That shifts every row below it up
# by one, so the previously-empty placeholder rows 7/8 inherit the
# (now-vacated) row numbers 6/7 and the old row 9's formatting, and the
# very last placeholder row (12) disappears entirely. The "#" column (A)
# is renumbered back to a clean 1..9 sequence, and the now-unused risk
# descriptions/strings fall out of the shared string table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Riscos")

# Delete the whole row - shifts rows 7..12 up into 6..11.
$ws.Rows.Item(6).Delete()

# The rows that shifted up (now 6, 7, 8) used to be blank placeholder
# rows further down the sheet; clear any data/formulas they picked up
# from their old position, keeping just the formatting that came along.
$ws.Range("B6:J8").ClearContents()

# Renumber the "#" column so it again reads 1..9 top to bottom.
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

# Restore the view: zoomed out a bit further, selection parked on D23.
$excel.ActiveWindow.Zoom = 70
$ws.Range("D23").Select()
